$d = $word.ActiveDocument

# Locate the target sentence in the last bullet paragraph.
$old = "Average windspeed around 10mphs across the globe"
$prefix = "Average windspeed "         # becomes its own run
$middle = "is "                        # new run inserted in the middle
$suffix = "around 10mphs across the globe"  # remaining original text, its own run

$full = $d.Content
$full.Find.Execute($old) | Out-Null
$start = $full.Start

$splitAt    = $start + $prefix.Length              # boundary between run 1 and run 2
$afterIns   = $splitAt + $middle.Length             # boundary between run 2 (new) and run 3

# The existing "_GoBack" bookmark currently sits at the very end of this
# sentence; remove it so it can be re-created at its new position once the
# text has been edited.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Insert the new word "is " right after "Average windspeed ".
$insertPoint = $d.Range($splitAt, $splitAt)
$insertPoint.InsertAfter($middle)

# Force a run boundary between "Average windspeed " and "is " by briefly
# dropping (and immediately removing) a bookmark at that position - the
# engine splits the underlying run when a bookmark is added/removed at an
# interior offset, and the split persists after the marker bookmark goes away.
$tempRange = $d.Range($splitAt, $splitAt)
$d.Bookmarks.Add("zzTempSplitMarker", $tempRange)
$d.Bookmarks("zzTempSplitMarker").Delete()

# Re-create "_GoBack" between the new "is " run and the remaining text - this
# both restores the bookmark and splits "is " away from "around 10mphs...".
$newBookmarkRange = $d.Range($afterIns, $afterIns)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)

Write-Output $d.Content.Text
